# SCD0174 fix: renumber CIF id column (A), rename customer placeholders
# (column B), and refresh the sheet's stored selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: CIF numbers -------------------------------------------------
# A2 is a plain literal; A3 becomes "=A2+1"; A4:A8 becomes a shared formula
# "=A3+1" filled down (matches the ref="A4:A8" si="1" group in the target).
$ws.Range("A2").Value = 9020304198
$ws.Range("A3").Formula = "=A2+1"
$ws.Range("A4:A8").Formula = "=A3+1"

# --- Column B: customer name placeholders ----------------------------------
$ws.Range("B2").Value = "dedic 36"
$ws.Range("B3").Value = "dedic 37"
$ws.Range("B4").Value = "dedic 38"
$ws.Range("B5").Value = "dedic 39"
$ws.Range("B6").Value = "dedic 40"
$ws.Range("B7").Value = "dedic 41"
$ws.Range("B8").Value = "dedic 42"

# --- Sheet view: move the saved selection to D10 ---------------------------
$ws.Range("D10").Select()
